$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.033443141075383
$ws.Range("D2").Value = 1.037770101111458
$ws.Range("E2").Value = 1.037085900938381
$ws.Range("F2").Value = 1.032078606655924
$ws.Range("I2").Value = 1.038069521806609
$ws.Range("J2").Value = 1.038567839859251
$ws.Range("K2").Value = 1.040560038278691
$ws.Range("L2").Value = 1.039877791016656
$ws.Range("M2").Value = 1.034884874677332
$ws.Range("N2").Value = 1.040042725317279
# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.034495350735429
$ws.Range("D3").Value = 1.038588155933879
$ws.Range("E3").Value = 1.038086392403691
$ws.Range("F3").Value = 1.03376161757437
$ws.Range("I3").Value = 1.038391595922896
$ws.Range("J3").Value = 1.039262244604023
$ws.Range("K3").Value = 1.041187947440163
$ws.Range("L3").Value = 1.040687512215347
$ws.Range("M3").Value = 1.036374244038986
$ws.Range("N3").Value = 1.040738116196406
# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.035175733481377
$ws.Range("D4").Value = 1.039116987653382
$ws.Range("E4").Value = 1.038733661922403
$ws.Range("F4").Value = 1.034850085502957
$ws.Range("I4").Value = 1.038598374468225
$ws.Range("J4").Value = 1.039710542792559
$ws.Range("K4").Value = 1.041593095281353
$ws.Range("L4").Value = 1.041210733069864
$ws.Range("M4").Value = 1.037336961699244
$ws.Range("N4").Value = 1.041187051019791
# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.035461656154014
$ws.Range("D5").Value = 1.039339188411004
$ws.Range("E5").Value = 1.039005747541528
$ws.Range("F5").Value = 1.035307552671882
$ws.Range("I5").Value = 1.038684915700268
$ws.Range("J5").Value = 1.039898762343401
$ws.Range("K5").Value = 1.041763144640651
$ws.Range("L5").Value = 1.04143052318062
$ws.Range("M5").Value = 1.037741454470166
$ws.Range("N5").Value = 1.041375537863984
# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.035509657411936
$ws.Range("D6").Value = 1.039376489872652
$ws.Range("E6").Value = 1.039051430393727
$ws.Range("F6").Value = 1.035384356326528
$ws.Range("I6").Value = 1.038699423589747
$ws.Range("J6").Value = 1.039930350899374
$ws.Range("K6").Value = 1.041791680586817
$ws.Range("L6").Value = 1.041467416857537
$ws.Range("M6").Value = 1.037809357071515
$ws.Range("N6").Value = 1.04140717127933
# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.035179554425245
$ws.Range("D7").Value = 1.039119957181869
$ws.Range("E7").Value = 1.038737297646444
$ws.Range("F7").Value = 1.034856198682454
$ws.Range("I7").Value = 1.038599532362149
$ws.Range("J7").Value = 1.039713058752738
$ws.Range("K7").Value = 1.041595368567927
$ws.Range("L7").Value = 1.041213670590285
$ws.Range("M7").Value = 1.037342367461677
$ws.Range("N7").Value = 1.041189570552921
# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.033798837990086
$ws.Range("D8").Value = 1.038046671425296
$ws.Range("E8").Value = 1.037424046283184
$ws.Range("F8").Value = 1.032647507069905
$ws.Range("I8").Value = 1.038178705048055
$ws.Range("J8").Value = 1.038802730751661
$ws.Range("K8").Value = 1.040772482032075
$ws.Range("L8").Value = 1.040151590451215
$ws.Range("M8").Value = 1.035388425873213
$ws.Range("N8").Value = 1.040277949781667
# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.031362186837383
$ws.Range("D9").Value = 1.036151510363911
$ws.Range("E9").Value = 1.035108994569903
$ws.Range("F9").Value = 1.028750938340536
$ws.Range("I9").Value = 1.037424683474921
$ws.Range("J9").Value = 1.037190691308793
$ws.Range("K9").Value = 1.039313593974527
$ws.Range("L9").Value = 1.038274481664023
$ws.Range("M9").Value = 1.031937346251021
$ws.Range("N9").Value = 1.038663621057887
# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.029735191198113
$ws.Range("D10").Value = 1.034885404923
$ws.Range("E10").Value = 1.033564923033724
$ws.Range("F10").Value = 1.026149676863493
$ws.Range("I10").Value = 1.036913581740449
$ws.Range("J10").Value = 1.036110596142776
$ws.Range("K10").Value = 1.038334990546673
$ws.Range("L10").Value = 1.037019244471553
$ws.Range("M10").Value = 1.029630842641893
$ws.Range("N10").Value = 1.037581992032846
# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.029030052121346
$ws.Range("D11").Value = 1.034336523977837
$ws.Range("E11").Value = 1.032896138443653
$ws.Range("F11").Value = 1.025022354140596
$ws.Range("I11").Value = 1.036690263126167
$ws.Range("J11").Value = 1.035641604533457
$ws.Range("K11").Value = 1.03790980529021
$ws.Range("L11").Value = 1.036474786798826
$ws.Range("M11").Value = 1.028630636498159
$ws.Range("N11").Value = 1.037112334401647
# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.028768033824035
$ws.Range("D12").Value = 1.03413254658251
$ws.Range("E12").Value = 1.032647691683701
$ws.Range("F12").Value = 1.024603462970263
$ws.Range("I12").Value = 1.036607010122237
$ws.Range("J12").Value = 1.035467202796281
$ws.Range("K12").Value = 1.037751654438174
$ws.Range("L12").Value = 1.036272409407586
$ws.Range("M12").Value = 1.028258885765503
$ws.Range("N12").Value = 1.036937684994
# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.028824242119842
$ws.Range("D13").Value = 1.034176304904355
$ws.Range("E13").Value = 1.032700985758375
$ws.Range("F13").Value = 1.024693323634251
$ws.Range("I13").Value = 1.036624881877424
$ws.Range("J13").Value = 1.035504621517751
$ws.Range("K13").Value = 1.037785588223141
$ws.Range("L13").Value = 1.036315826477074
$ws.Range("M13").Value = 1.028338638080971
$ws.Range("N13").Value = 1.036975156854347
# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.029008395624223
$ws.Range("D14").Value = 1.034319665159867
$ws.Range("E14").Value = 1.032875602386291
$ws.Range("F14").Value = 1.024987731646735
$ws.Range("I14").Value = 1.036683387583123
$ws.Range("J14").Value = 1.035627192466905
$ws.Range("K14").Value = 1.0378967369482
$ws.Range("L14").Value = 1.036458061117786
$ws.Range("M14").Value = 1.028599912197911
$ws.Range("N14").Value = 1.037097901868307
# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.029121845590908
$ws.Range("D15").Value = 1.034407981050328
$ws.Range("E15").Value = 1.032983185378532
$ws.Range("F15").Value = 1.025169105548184
$ws.Range("I15").Value = 1.036719394766891
$ws.Range("J15").Value = 1.035702686273198
$ws.Range("K15").Value = 1.037965190406085
$ws.Range("L15").Value = 1.036545677766973
$ws.Range("M15").Value = 1.028760861133629
$ws.Range("N15").Value = 1.037173502884465
# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.029781975300137
$ws.Range("D16").Value = 1.034921818579643
$ws.Range("E16").Value = 1.033609303892835
$ws.Range("F16").Value = 1.026224472564731
$ws.Range("I16").Value = 1.036928360281361
$ws.Range("J16").Value = 1.036141693979159
$ws.Range("K16").Value = 1.038363178172694
$ws.Range("L16").Value = 1.037055358568867
$ws.Range("M16").Value = 1.029697191269401
$ws.Range("N16").Value = 1.037613134031723
# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.030195884531358
$ws.Range("D17").Value = 1.035243960792249
$ws.Range("E17").Value = 1.034001999102647
$ws.Range("F17").Value = 1.026886213102772
$ws.Range("I17").Value = 1.037058900602357
$ws.Range("J17").Value = 1.036416721887937
$ws.Range("K17").Value = 1.038612438040007
$ws.Range("L17").Value = 1.037374817311271
$ws.Range("M17").Value = 1.030284125893637
$ws.Range("N17").Value = 1.037888552511678
# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.030437249256041
$ws.Range("D18").Value = 1.035431798355089
$ws.Range("E18").Value = 1.034231033214826
$ws.Range("F18").Value = 1.027272103219706
$ws.Range("I18").Value = 1.03713484885713
$ws.Range("J18").Value = 1.036577015372905
$ws.Range("K18").Value = 1.038757687909456
$ws.Range("L18").Value = 1.037561062404559
$ws.Range("M18").Value = 1.030626333169913
$ws.Range("N18").Value = 1.03804907363178
# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.030519538020688
$ws.Range("D19").Value = 1.035495835534633
$ws.Range("E19").Value = 1.034309124863377
$ws.Range("F19").Value = 1.027403666460838
$ws.Range("I19").Value = 1.037160712418263
$ws.Range("J19").Value = 1.036631650045542
$ws.Range("K19").Value = 1.038807190789294
$ws.Range("L19").Value = 1.037624551987324
$ws.Range("M19").Value = 1.030742993207136
$ws.Range("N19").Value = 1.038103785891919
# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.030151482367723
$ws.Range("D20").Value = 1.035209404457869
$ws.Range("E20").Value = 1.033959868533564
$ws.Range("F20").Value = 1.026815224181161
$ws.Range("I20").Value = 1.037044914903013
$ws.Range("J20").Value = 1.036387227002892
$ws.Range("K20").Value = 1.038585709221467
$ws.Range("L20").Value = 1.037340551709808
$ws.Range("M20").Value = 1.030221168083334
$ws.Range("N20").Value = 1.037859015740512
# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.028954169737044
$ws.Range("D21").Value = 1.034277451867622
$ws.Range("E21").Value = 1.032824182993502
$ws.Range("F21").Value = 1.024901040131331
$ws.Range("I21").Value = 1.036666167468863
$ws.Range("J21").Value = 1.035591103838473
$ws.Range("K21").Value = 1.037864012450119
$ws.Range("L21").Value = 1.036416180482093
$ws.Range("M21").Value = 1.028522979886578
$ws.Range("N21").Value = 1.037061761989883
# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.028200802227122
$ws.Range("D22").Value = 1.033690925424945
$ws.Range("E22").Value = 1.03210995580103
$ws.Range("F22").Value = 1.023696625012852
$ws.Range("I22").Value = 1.03642628318367
$ws.Range("J22").Value = 1.035089406888714
$ws.Range("K22").Value = 1.037408990465987
$ws.Range("L22").Value = 1.035834171724871
$ws.Range("M22").Value = 1.027453928936851
$ws.Range("N22").Value = 1.03655935257291
# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.028600231326533
$ws.Range("D23").Value = 1.034001908664863
$ws.Range("E23").Value = 1.032488598438143
$ws.Range("F23").Value = 1.024335195854719
$ws.Range("I23").Value = 1.036553616607716
$ws.Range("J23").Value = 1.035355474778511
$ws.Range("K23").Value = 1.037650326326263
$ws.Range("L23").Value = 1.036142783830502
$ws.Range("M23").Value = 1.028020782093317
$ws.Range("N23").Value = 1.036825798309629
# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.030171545983712
$ws.Range("D24").Value = 1.035225019170788
$ws.Range("E24").Value = 1.033978905580492
$ws.Range("F24").Value = 1.026847301301135
$ws.Range("I24").Value = 1.037051235036801
$ws.Range("J24").Value = 1.036400554859612
$ws.Range("K24").Value = 1.038597787253741
$ws.Range("L24").Value = 1.037356035137125
$ws.Range("M24").Value = 1.030249616443963
$ws.Range("N24").Value = 1.037872362524318
# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.031992563245986
$ws.Range("D25").Value = 1.036641921153986
$ws.Range("E25").Value = 1.035707609411423
$ws.Range("F25").Value = 1.029758885322443
$ws.Range("I25").Value = 1.037621097627511
$ws.Range("J25").Value = 1.037608388870504
$ws.Range("K25").Value = 1.03969180686734
$ws.Range("L25").Value = 1.038760429080699
$ws.Range("M25").Value = 1.03283052193367
$ws.Range("N25").Value = 1.039081911798045
